# Auto-generated script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.017.68'
$ws.Range('E2').Value = '  -0.74%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.222.70'
$ws.Range('E3').Value = '  -1.31%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.07'
$ws.Range('E5').Value = '  -1.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').Value = '  -1.54%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.39'
$ws.Range('E7').Value = '  +0.28%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  +0.04%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.32'
$ws.Range('E10').Value = '  +5.68%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0964'
$ws.Range('E11').Value = '  +2.10%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.19'
$ws.Range('E12').Value = '  +0.50%  '

$ws.Range('E13').Value = '  +0.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.554.65'
$ws.Range('E14').Value = '  -1.31%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.32'
$ws.Range('E15').Value = '  -0.97%  '

$ws.Range('E16').Value = '  -0.99%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.215.78'
$ws.Range('E17').Value = '  -1.72%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.955.67'
$ws.Range('E18').Value = '  -0.53%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('E19').Value = '  +12.18%  '

$ws.Range('E20').Value = '  +1.56%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.57'
$ws.Range('E21').Value = '  +0.85%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.28'
$ws.Range('E22').Value = '  +38.34%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.27'
$ws.Range('E23').Value = '  +0.01%  '

$ws.Range('E24').Value = '  -7.70%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.62'
$ws.Range('E25').Value = '  +4.59%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.59'
$ws.Range('E27').Value = '  +0.58%  '

$ws.Range('E28').Value = '  -1.38%  '

$ws.Range('E29').Value = '  -3.45%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.78'
$ws.Range('E30').Value = '  -1.56%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.65'
$ws.Range('E31').Value = '  -0.27%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.76'
$ws.Range('E32').Value = '  +16.85%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0804'
$ws.Range('E33').Value = '  -2.46%  '

$ws.Range('E34').Value = '  -0.15%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.53'
$ws.Range('E35').Value = '  -3.90%  '

$ws.Range('E36').Value = '  -4.28%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.30'
$ws.Range('E37').Value = '  -5.36%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0302'
$ws.Range('E38').Value = '  -1.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.02'
$ws.Range('E39').Value = '  -3.69%  '

$ws.Range('E40').Value = '  -1.72%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '65.15'
$ws.Range('E41').Value = '  +5.88%  '

$ws.Range('E42').Value = '  -2.16%  '

$ws.Range('E43').Value = '  -1.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.77'
$ws.Range('E44').Value = '  +1.25%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.55'
$ws.Range('E45').Value = '  -2.95%  '

$ws.Range('E46').Value = '  +0.22%  '

$ws.Range('E47').Value = '  +7.38%  '

$ws.Range('E48').Value = '  -0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.17'
$ws.Range('E49').Value = '  -0.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.71'
$ws.Range('E50').Value = '  +0.59%  '

$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.03'
$ws.Range('E51').Value = '  -2.26%  '

